# Generate Report for Archive
#
# The localization status "Ready for handoff" has moved on to
# "In Translation" for every tracked file: update the status values on
# the Overview roll-up sheet (zh-cn column E, de-de column F) and on the
# per-locale detail sheets (zh-cn, de-de - Status is column C), then
# shrink the now-narrower "Status" columns to their new auto-fit width.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn (col E) / de-de (col F) status columns ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F4").Value = "In Translation"
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn detail sheet: Status column (col C) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C4").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de detail sheet: Status column (col C) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C4").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.5
